# Refresh LR-pair (Gnai2-Oprd1) NATMI TPM output with new expression values.
# Rows 2-6 get updated figures for the existing Sending-cluster -> FAPs pairs
# (figures shift: the row that used to belong to each sending cluster now
# carries the next cluster's numbers), and rows 7-11 are newly added
# Sending-cluster -> Resolving-Mac pairs (mirrored MuSCs/FAPs row duplicated).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 197.5433703333333
$ws.Cells.Item(2, 8).Value = 592.6301109999999
$ws.Cells.Item(2, 9).Value = 0.3388703761585983
$ws.Cells.Item(2, 10).Value = 0.3388703761585982
$ws.Cells.Item(2, 15).Value = 0.7290280598220596
$ws.Cells.Item(2, 16).Value = 0.7290280598220598
$ws.Cells.Item(2, 17).Value = 4.278723553629888
$ws.Cells.Item(2, 18).Value = 38.50851198266899
$ws.Cells.Item(2, 19).Value = 0.2470460128620744
$ws.Cells.Item(2, 20).Value = 0.2470460128620744
# Row 3
$ws.Cells.Item(3, 1).Value = 'ECs'
$ws.Cells.Item(3, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(3, 7).Value = 197.5433703333333
$ws.Cells.Item(3, 8).Value = 592.6301109999999
$ws.Cells.Item(3, 9).Value = 0.3388703761585983
$ws.Cells.Item(3, 10).Value = 0.3388703761585982
$ws.Cells.Item(3, 13).Value = 0.008050666666666666
$ws.Cells.Item(3, 14).Value = 0.024152
$ws.Cells.Item(3, 15).Value = 0.2709719401779404
$ws.Cells.Item(3, 16).Value = 0.2709719401779404
$ws.Cells.Item(3, 17).Value = 1.590355826763555
$ws.Cells.Item(3, 18).Value = 14.313202440872
$ws.Cells.Item(3, 19).Value = 0.09182436329652384
$ws.Cells.Item(3, 20).Value = 0.09182436329652383
# Row 4
$ws.Cells.Item(4, 1).Value = 'FAPs'
$ws.Cells.Item(4, 7).Value = 79.82725266666667
$ws.Cells.Item(4, 8).Value = 239.481758
$ws.Cells.Item(4, 9).Value = 0.1369374790620155
$ws.Cells.Item(4, 10).Value = 0.1369374790620154
$ws.Cells.Item(4, 15).Value = 0.7290280598220596
$ws.Cells.Item(4, 16).Value = 0.7290280598220598
$ws.Cells.Item(4, 17).Value = 1.729031683675778
$ws.Cells.Item(4, 18).Value = 15.561285153082
$ws.Cells.Item(4, 19).Value = 0.09983126467750505
$ws.Cells.Item(4, 20).Value = 0.09983126467750505
# Row 5
$ws.Cells.Item(5, 1).Value = 'FAPs'
$ws.Cells.Item(5, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(5, 7).Value = 79.82725266666667
$ws.Cells.Item(5, 8).Value = 239.481758
$ws.Cells.Item(5, 9).Value = 0.1369374790620155
$ws.Cells.Item(5, 10).Value = 0.1369374790620154
$ws.Cells.Item(5, 13).Value = 0.008050666666666666
$ws.Cells.Item(5, 14).Value = 0.024152
$ws.Cells.Item(5, 15).Value = 0.2709719401779404
$ws.Cells.Item(5, 16).Value = 0.2709719401779404
$ws.Cells.Item(5, 17).Value = 0.6426626021351111
$ws.Cells.Item(5, 18).Value = 5.783963419216
$ws.Cells.Item(5, 19).Value = 0.03710621438451042
$ws.Cells.Item(5, 20).Value = 0.03710621438451041
# Row 6
$ws.Cells.Item(6, 1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(6, 7).Value = 148.824417
$ws.Cells.Item(6, 8).Value = 446.473251
$ws.Cells.Item(6, 9).Value = 0.2552967790580629
$ws.Cells.Item(6, 10).Value = 0.2552967790580629
$ws.Cells.Item(6, 15).Value = 0.7290280598220596
$ws.Cells.Item(6, 16).Value = 0.7290280598220598
$ws.Cells.Item(6, 17).Value = 3.223487264081
$ws.Cells.Item(6, 18).Value = 29.011385376729
$ws.Cells.Item(6, 19).Value = 0.1861185155155206
$ws.Cells.Item(6, 20).Value = 0.1861185155155206
# Row 7
$ws.Cells.Item(7, 1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(7, 2).Value = 'Gnai2'
$ws.Cells.Item(7, 3).Value = 'Oprd1'
$ws.Cells.Item(7, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 148.824417
$ws.Cells.Item(7, 8).Value = 446.473251
$ws.Cells.Item(7, 9).Value = 0.2552967790580629
$ws.Cells.Item(7, 10).Value = 0.2552967790580629
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.008050666666666666
$ws.Cells.Item(7, 14).Value = 0.024152
$ws.Cells.Item(7, 15).Value = 0.2709719401779404
$ws.Cells.Item(7, 16).Value = 0.2709719401779404
$ws.Cells.Item(7, 17).Value = 1.198135773128
$ws.Cells.Item(7, 18).Value = 10.783221958152
$ws.Cells.Item(7, 19).Value = 0.06917826354254226
$ws.Cells.Item(7, 20).Value = 0.06917826354254226
# Row 8
$ws.Cells.Item(8, 1).Value = 'MuSCs'
$ws.Cells.Item(8, 2).Value = 'Gnai2'
$ws.Cells.Item(8, 3).Value = 'Oprd1'
$ws.Cells.Item(8, 4).Value = 'FAPs'
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 35.426853
$ws.Cells.Item(8, 8).Value = 106.280559
$ws.Cells.Item(8, 9).Value = 0.06077202683121193
$ws.Cells.Item(8, 10).Value = 0.06077202683121192
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.02165966666666666
$ws.Cells.Item(8, 14).Value = 0.064979
$ws.Cells.Item(8, 15).Value = 0.7290280598220596
$ws.Cells.Item(8, 16).Value = 0.7290280598220598
$ws.Cells.Item(8, 17).Value = 0.7673338270289999
$ws.Cells.Item(8, 18).Value = 6.906004443260999
$ws.Cells.Item(8, 19).Value = 0.04430451281221258
$ws.Cells.Item(8, 20).Value = 0.04430451281221259
# Row 9
$ws.Cells.Item(9, 1).Value = 'MuSCs'
$ws.Cells.Item(9, 2).Value = 'Gnai2'
$ws.Cells.Item(9, 3).Value = 'Oprd1'
$ws.Cells.Item(9, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 35.426853
$ws.Cells.Item(9, 8).Value = 106.280559
$ws.Cells.Item(9, 9).Value = 0.06077202683121193
$ws.Cells.Item(9, 10).Value = 0.06077202683121192
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.008050666666666666
$ws.Cells.Item(9, 14).Value = 0.024152
$ws.Cells.Item(9, 15).Value = 0.2709719401779404
$ws.Cells.Item(9, 16).Value = 0.2709719401779404
$ws.Cells.Item(9, 17).Value = 0.285209784552
$ws.Cells.Item(9, 18).Value = 2.566888060968
$ws.Cells.Item(9, 19).Value = 0.01646751401899934
$ws.Cells.Item(9, 20).Value = 0.01646751401899934
# Row 10
$ws.Cells.Item(10, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(10, 2).Value = 'Gnai2'
$ws.Cells.Item(10, 3).Value = 'Oprd1'
$ws.Cells.Item(10, 4).Value = 'FAPs'
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 121.3248153333333
$ws.Cells.Item(10, 8).Value = 363.974446
$ws.Cells.Item(10, 9).Value = 0.2081233388901116
$ws.Cells.Item(10, 10).Value = 0.2081233388901115
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.02165966666666666
$ws.Cells.Item(10, 14).Value = 0.064979
$ws.Cells.Item(10, 15).Value = 0.7290280598220596
$ws.Cells.Item(10, 16).Value = 0.7290280598220598
$ws.Cells.Item(10, 17).Value = 2.627855058514889
$ws.Cells.Item(10, 18).Value = 23.650695526634
$ws.Cells.Item(10, 19).Value = 0.1517277539547471
$ws.Cells.Item(10, 20).Value = 0.1517277539547471
# Row 11
$ws.Cells.Item(11, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(11, 2).Value = 'Gnai2'
$ws.Cells.Item(11, 3).Value = 'Oprd1'
$ws.Cells.Item(11, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 121.3248153333333
$ws.Cells.Item(11, 8).Value = 363.974446
$ws.Cells.Item(11, 9).Value = 0.2081233388901116
$ws.Cells.Item(11, 10).Value = 0.2081233388901115
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.008050666666666666
$ws.Cells.Item(11, 14).Value = 0.024152
$ws.Cells.Item(11, 15).Value = 0.2709719401779404
$ws.Cells.Item(11, 16).Value = 0.2709719401779404
$ws.Cells.Item(11, 17).Value = 0.9767456466435555
$ws.Cells.Item(11, 18).Value = 8.790710819792
$ws.Cells.Item(11, 19).Value = 0.05639558493536451
$ws.Cells.Item(11, 20).Value = 0.05639558493536451
